$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-5) are cyclically rotated:
#   new row2 <- old row4
#   new row3 <- old row2
#   new row4 <- old row5
#   new row5 <- old row3
# Only columns D (Fecha), J (Volumen), K (Precio minimo), L (Precio maximo),
# M (Precio promedio ponderado) and P (Precio $/Kg) actually change value.

$cols = @("D", "J", "K", "L", "M", "P")

# Capture the "before" values for each relevant column/row first since we will
# overwrite them in place.
$orig = @{}
foreach ($col in $cols) {
    $orig[$col] = @{}
    for ($r = 2; $r -le 5; $r++) {
        $orig[$col][$r] = $ws.Range("$col$r").Value2
    }
}

# Row -> source row mapping (new row gets the old values of the source row)
$rowMap = @{ 2 = 4; 3 = 2; 4 = 5; 5 = 3 }

foreach ($col in $cols) {
    foreach ($destRow in 2..5) {
        $srcRow = $rowMap[$destRow]
        $ws.Range("$col$destRow").Value2 = $orig[$col][$srcRow]
    }
}
